$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 15.67790914277992
$ws.Cells.Item(3, 1).Value = 19.88800271446007
$ws.Cells.Item(4, 1).Value = 18.89080234373404
$ws.Cells.Item(5, 1).Value = 20.28454590844126
$ws.Cells.Item(6, 1).Value = 10.55267256212903
$ws.Cells.Item(7, 1).Value = 20.63310720507212
$ws.Cells.Item(8, 1).Value = 16.58696707729638
$ws.Cells.Item(9, 1).Value = 15.38926229123268
$ws.Cells.Item(10, 1).Value = 21.34770284898741
$ws.Cells.Item(11, 1).Value = 22.6759862887024
$ws.Cells.Item(12, 1).Value = 13.13312222896269
$ws.Cells.Item(13, 1).Value = 18.45447784467422
$ws.Cells.Item(14, 1).Value = 16.64593431388198
$ws.Cells.Item(15, 1).Value = 8.17536073061018
$ws.Cells.Item(16, 1).Value = 8.947183266864215
$ws.Cells.Item(17, 1).Value = 4.398575093383272
$ws.Cells.Item(18, 1).Value = 4.398575093383272
$ws.Cells.Item(19, 1).Value = 4.398575093383272
$ws.Cells.Item(20, 1).Value = 4.398575093383272
$ws.Cells.Item(21, 1).Value = 4.398575093383272
$ws.Cells.Item(22, 1).Value = 4.398575093383272
$ws.Cells.Item(23, 1).Value = 4.398575093383272
$ws.Cells.Item(24, 1).Value = 4.398575093383272
$ws.Cells.Item(25, 1).Value = 4.398575093383272
$ws.Cells.Item(26, 1).Value = 4.398575093383272
$ws.Cells.Item(27, 1).Value = 4.398575093383272
$ws.Cells.Item(28, 1).Value = 4.398575093383272
$ws.Cells.Item(29, 1).Value = 4.398575093383272
$ws.Cells.Item(30, 1).Value = 4.398575093383272
$ws.Cells.Item(31, 1).Value = 4.398575093383272
$ws.Cells.Item(32, 1).Value = 4.398575093383272
$ws.Cells.Item(33, 1).Value = 4.398575093383272
$ws.Cells.Item(34, 1).Value = 4.398575093383272
$ws.Cells.Item(35, 1).Value = 4.398575093383272
$ws.Cells.Item(36, 1).Value = 4.398575093383272
$ws.Cells.Item(37, 1).Value = 4.398575093383272
$ws.Cells.Item(38, 1).Value = 4.398575093383272
$ws.Cells.Item(39, 1).Value = 4.398575093383272
$ws.Cells.Item(40, 1).Value = 4.398575093383272
$ws.Cells.Item(41, 1).Value = 4.398575093383272
$ws.Cells.Item(42, 1).Value = 4.398575093383272
$ws.Cells.Item(43, 1).Value = 4.398575093383272
$ws.Cells.Item(44, 1).Value = 4.398575093383272
$ws.Cells.Item(45, 1).Value = 4.398575093383272
$ws.Cells.Item(46, 1).Value = 4.398575093383272
$ws.Cells.Item(47, 1).Value = 4.398575093383272
$ws.Cells.Item(48, 1).Value = 4.398575093383272
$ws.Cells.Item(49, 1).Value = 4.398575093383272
$ws.Cells.Item(50, 1).Value = 4.398575093383272
$ws.Cells.Item(51, 1).Value = 4.398575093383272
$ws.Cells.Item(52, 1).Value = 4.398575093383272
$ws.Cells.Item(53, 1).Value = 4.398575093383272
$ws.Cells.Item(54, 1).Value = 4.398575093383272
$ws.Cells.Item(55, 1).Value = 4.398575093383272
$ws.Cells.Item(56, 1).Value = 4.398575093383272
$ws.Cells.Item(57, 1).Value = 4.398575093383272
$ws.Cells.Item(58, 1).Value = 4.398575093383272
$ws.Cells.Item(59, 1).Value = 4.398575093383272
$ws.Cells.Item(60, 1).Value = 4.398575093383272
$ws.Cells.Item(61, 1).Value = 4.398575093383272
$ws.Cells.Item(62, 1).Value = 4.398575093383272
$ws.Cells.Item(63, 1).Value = 4.398575093383272
$ws.Cells.Item(64, 1).Value = 4.398575093383272
$ws.Cells.Item(65, 1).Value = 4.398575093383272
$ws.Cells.Item(66, 1).Value = 4.398575093383272
$ws.Cells.Item(67, 1).Value = 4.398575093383272
$ws.Cells.Item(68, 1).Value = 4.398575093383272
$ws.Cells.Item(69, 1).Value = 4.398575093383272
$ws.Cells.Item(70, 1).Value = 4.398575093383272
$ws.Cells.Item(71, 1).Value = 4.398575093383272
$ws.Cells.Item(72, 1).Value = 4.398575093383272
$ws.Cells.Item(73, 1).Value = 4.398575093383272
$ws.Cells.Item(74, 1).Value = 4.398575093383272
$ws.Cells.Item(75, 1).Value = 4.398575093383272
$ws.Cells.Item(76, 1).Value = 4.398575093383272
$ws.Cells.Item(77, 1).Value = 4.398575093383272
$ws.Cells.Item(78, 1).Value = 4.398575093383272
$ws.Cells.Item(79, 1).Value = 4.398575093383272
$ws.Cells.Item(80, 1).Value = 4.398575093383272
$ws.Cells.Item(81, 1).Value = 4.398575093383272
$ws.Cells.Item(82, 1).Value = 4.398575093383272
$ws.Cells.Item(83, 1).Value = 4.398575093383272
$ws.Cells.Item(84, 1).Value = 4.398575093383272
$ws.Cells.Item(85, 1).Value = 4.398575093383272
$ws.Cells.Item(86, 1).Value = 4.398575093383272
